# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos sheet with
# the latest scraped snapshot. Values are stored as plain text in the sheet
# (e.g. "321.30", "5.35%"), so each write uses a leading apostrophe to force
# text interpretation (otherwise Excel would coerce "321.30" into the number
# 321.3 and drop the trailing zero / exact formatting). The Style reset
# afterwards clears the quote-prefix formatting flag that the apostrophe
# trick leaves behind, so the cell's style stays identical to before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "321.30"
Set-TextValue 2 5 "5.35%"
Set-TextValue 3 4 "36.23"
Set-TextValue 3 5 "0.46%"
Set-TextValue 4 4 "5.122"
Set-TextValue 4 5 "1.83%"
Set-TextValue 5 4 "0.08095"
Set-TextValue 5 5 "3.03%"
Set-TextValue 6 4 "2.157"
Set-TextValue 6 5 "-0.47%"
Set-TextValue 7 4 "8.010"
Set-TextValue 7 5 "1.27%"
Set-TextValue 8 4 "0.9285"
Set-TextValue 8 5 "1.05%"
Set-TextValue 9 4 "0.1004"
Set-TextValue 9 5 "3.01%"
Set-TextValue 10 4 "0.1893"
Set-TextValue 10 5 "1.70%"
Set-TextValue 11 4 "0.09195"
Set-TextValue 11 5 "5.57%"
Set-TextValue 12 4 "0.03593"
Set-TextValue 12 5 "3.21%"
Set-TextValue 13 4 "0.09939"
Set-TextValue 13 5 "0.29%"
Set-TextValue 14 4 "0.001436"
Set-TextValue 14 5 "-0.16%"
Set-TextValue 15 4 "0.005682"
Set-TextValue 15 5 "-0.20%"
Set-TextValue 16 4 "3.446"
Set-TextValue 16 5 "-0.51%"
Set-TextValue 17 4 "4.140"
Set-TextValue 17 5 "1.46%"
Set-TextValue 18 4 "2.797"
Set-TextValue 18 5 "17.65%"
Set-TextValue 19 5 "-1.50%"
Set-TextValue 20 4 "0.1330"
Set-TextValue 20 5 "-1.05%"
Set-TextValue 21 4 "5.056"
Set-TextValue 21 5 "5.98%"
Set-TextValue 22 4 "0.2204"
Set-TextValue 22 5 "-0.04%"
Set-TextValue 23 4 "0.04600"
Set-TextValue 23 5 "-0.06%"
Set-TextValue 24 4 "0.001242"
Set-TextValue 24 5 "0.71%"
Set-TextValue 25 4 "0.004736"
Set-TextValue 25 5 "-7.05%"
Set-TextValue 26 5 "-7.05%"
Set-TextValue 27 4 "0.0004503"
Set-TextValue 27 5 "-5.16%"
Set-TextValue 39 4 "0.02027"
Set-TextValue 39 5 "10.76%"
Set-TextValue 40 4 "0.04989"
Set-TextValue 40 5 "4.59%"
Set-TextValue 41 4 "0.007809"
Set-TextValue 41 5 "1.08%"
Set-TextValue 42 4 "0.1401"
Set-TextValue 42 5 "0.29%"
Set-TextValue 43 4 "0.007820"
Set-TextValue 43 5 "1.14%"
Set-TextValue 44 4 "0.002082"
Set-TextValue 44 5 "-6.91%"
Set-TextValue 45 5 "8.96%"
Set-TextValue 46 4 "0.00006438"
Set-TextValue 46 5 "1.17%"
Set-TextValue 47 5 "0.11%"
Set-TextValue 48 5 "17.46%"
Set-TextValue 49 4 "0.001901"
Set-TextValue 49 5 "-4.89%"
Set-TextValue 50 4 "0.00002101"
Set-TextValue 50 5 "0.11%"
Set-TextValue 51 4 "0.0002001"
Set-TextValue 51 5 "0.11%"
